$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value2 = "'" + $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "20.563.60"
$ws.Range("E2").Value2 = "  +1.47%  "
Set-TextValue "D3" "1.476.05"
$ws.Range("E3").Value2 = "  +2.01%  "
Set-TextValue "D4" "1.010"
$ws.Range("E4").Value2 = "  +0.29%  "
Set-TextValue "D5" "0.9588"
$ws.Range("E5").Value2 = "  +2.72%  "
Set-TextValue "D6" "277.23"
$ws.Range("E6").Value2 = "  +1.20%  "
Set-TextValue "D7" "0.3536"
$ws.Range("E7").Value2 = "  -2.58%  "
Set-TextValue "D8" "0.3072"
$ws.Range("E8").Value2 = "  -0.11%  "
Set-TextValue "D9" "1.080"
$ws.Range("E9").Value2 = "  +5.39%  "
Set-TextValue "D10" "39.34"
$ws.Range("E10").Value2 = "  -0.66%  "
Set-TextValue "D11" "0.06642"
$ws.Range("E11").Value2 = "  +2.08%  "
Set-TextValue "D12" "1.005"
$ws.Range("E12").Value2 = "  +0.68%  "
$ws.Range("B13").Value2 = "Solana"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D13" "18.12"
$ws.Range("E13").Value2 = "  +3.41%  "
$ws.Range("B14").Value2 = "Polkadot"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "5.458"
$ws.Range("E14").Value2 = "  +2.06%  "
Set-TextValue "D15" "6.170"
$ws.Range("E15").Value2 = "  +1.78%  "
Set-TextValue "D16" "0.9584"
$ws.Range("E16").Value2 = "  +0.99%  "
Set-TextValue "D17" "0.00001013"
$ws.Range("E17").Value2 = "  +0.11%  "
Set-TextValue "D18" "1.475.55"
$ws.Range("E18").Value2 = "  +2.23%  "
Set-TextValue "D19" "0.05984"
$ws.Range("E19").Value2 = "  +5.63%  "
Set-TextValue "D20" "69.14"
$ws.Range("E20").Value2 = "  +0.11%  "
Set-TextValue "D21" "5.479"
$ws.Range("E21").Value2 = "  +1.60%  "
Set-TextValue "D22" "14.48"
$ws.Range("E22").Value2 = "  +1.54%  "
Set-TextValue "D23" "11.08"
$ws.Range("E23").Value2 = "  +2.79%  "
Set-TextValue "D24" "2.269"
$ws.Range("E24").Value2 = "  +0.55%  "
Set-TextValue "D25" "20.604.00"
$ws.Range("E25").Value2 = "  +1.54%  "
Set-TextValue "D26" "147.28"
$ws.Range("E26").Value2 = "  +5.00%  "
Set-TextValue "D27" "2.081"
$ws.Range("E27").Value2 = "  +1.12%  "
$ws.Range("E28").Value2 = "  +1.05%  "
Set-TextValue "D29" "1.638.28"
Set-TextValue "D30" "114.49"
$ws.Range("E30").Value2 = "  +3.38%  "
Set-TextValue "D31" "3.912"
$ws.Range("E31").Value2 = "  -1.53%  "
Set-TextValue "D32" "4.932"
$ws.Range("E32").Value2 = "  +2.62%  "
$ws.Range("B33").Value2 = "Stellar"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D33" "0.07910"
$ws.Range("E33").Value2 = "  +2.70%  "
$ws.Range("B34").Value2 = "ImmutableX"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D34" "0.7956"
$ws.Range("E34").Value2 = "  +1.51%  "
Set-TextValue "D35" "1.194"
$ws.Range("E35").Value2 = "  +7.42%  "
Set-TextValue "D36" "1.434"
$ws.Range("E36").Value2 = "  -1.54%  "
Set-TextValue "D37" "0.05672"
$ws.Range("E37").Value2 = "  +0.47%  "
Set-TextValue "D38" "4.695"
$ws.Range("E38").Value2 = "  +1.19%  "
Set-TextValue "D39" "0.9591"
$ws.Range("E39").Value2 = "  +2.05%  "
Set-TextValue "D40" "0.02015"
$ws.Range("E40").Value2 = "  +0.76%  "
Set-TextValue "D41" "10.22"
$ws.Range("E41").Value2 = "  +0.20%  "
$ws.Range("B42").Value2 = "Algorand"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D42" "0.1842"
$ws.Range("E42").Value2 = "  -0.13%  "
$ws.Range("B43").Value2 = "FraxShare"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "7.336"
$ws.Range("E43").Value2 = "  +4.93%  "
$ws.Range("B44").Value2 = "PancakeSwap"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D44" "3.509"
$ws.Range("E44").Value2 = "  +1.06%  "
$ws.Range("B45").Value2 = "TheSandbox"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D45" "0.5211"
$ws.Range("E45").Value2 = "  -0.05%  "
Set-TextValue "D46" "12.00"
$ws.Range("E46").Value2 = "  +1.79%  "
Set-TextValue "D47" "119.64"
$ws.Range("E47").Value2 = "  +2.40%  "
Set-TextValue "D48" "0.5153"
$ws.Range("E48").Value2 = "  +1.17%  "
Set-TextValue "D49" "1.806"
$ws.Range("E49").Value2 = "  +4.30%  "
Set-TextValue "D50" "0.06399"
$ws.Range("E50").Value2 = "  +0.43%  "
Set-TextValue "D51" "0.9916"
$ws.Range("E51").Value2 = "  +0.94%  "
